# Helper: pack an RGB triple into the COM BGR-in-int "RGB" long used by
# ColorFormat.RGB (0x00BBGGRR).
function ToCOMColor($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 6 table: switch the table style (Design > Table Styles gallery)
#    from the old custom "Table_0" style to the built-in style whose
#    GUID is {DE2031D8-34A4-4DAC-A049-9B5AE11A9E55}.
# ---------------------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{DE2031D8-34A4-4DAC-A049-9B5AE11A9E55}")

# ---------------------------------------------------------------------
# 2) Presentation theme: swap the custom "Integral" theme colors for the
#    stock Office theme palette (Design tab > Themes > Office Theme).
# ---------------------------------------------------------------------
$cs = $p.SlideMaster.ColorScheme
$cs.Colors(1).RGB  = ToCOMColor 0x00 0x00 0x00   # Dark 1    -> 000000
$cs.Colors(2).RGB  = ToCOMColor 0xFF 0xFF 0xFF   # Light 1   -> FFFFFF
$cs.Colors(3).RGB  = ToCOMColor 0x44 0x54 0x6A   # Dark 2    -> 44546A
$cs.Colors(4).RGB  = ToCOMColor 0xE7 0xE6 0xE6   # Light 2   -> E7E6E6
$cs.Colors(5).RGB  = ToCOMColor 0x5B 0x9B 0xD5   # Accent 1  -> 5B9BD5
$cs.Colors(6).RGB  = ToCOMColor 0xED 0x7D 0x31   # Accent 2  -> ED7D31
$cs.Colors(7).RGB  = ToCOMColor 0xA5 0xA5 0xA5   # Accent 3  -> A5A5A5
$cs.Colors(8).RGB  = ToCOMColor 0xFF 0xC0 0x00   # Accent 4  -> FFC000
$cs.Colors(9).RGB  = ToCOMColor 0x44 0x72 0xC4   # Accent 5  -> 4472C4
$cs.Colors(10).RGB = ToCOMColor 0x70 0xAD 0x47   # Accent 6  -> 70AD47
$cs.Colors(11).RGB = ToCOMColor 0x05 0x63 0xC1   # Hyperlink -> 0563C1
$cs.Colors(12).RGB = ToCOMColor 0x95 0x4F 0x72   # Followed Hyperlink -> 954F72
